# Sargatanas_Profits: refresh the per-item market snapshot + computed Leve profit
# columns (H:currentAveragePrice, I/J:NQ/HQ current average, K/L:Leve NQ/HQ price,
# M/N:Leve NQ/HQ profit) on each job sheet, as produced by the scheduled market-data
# runner.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 503750
$ws.Range("I76").Value = 503750
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 503750
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -503435
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 503750
$ws.Range("I79").Value = 503750
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 503750
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -502658
$ws.Range("N79").ClearContents()
$ws.Range("H116").Value = 16675886
$ws.Range("I116").Value = 62503972
$ws.Range("J116").Value = 11127.091
$ws.Range("K116").Value = 62503972
$ws.Range("L116").Value = 11127.091
$ws.Range("M116").Value = -62500530
$ws.Range("N116").Value = -18011.091
$ws.Range("H127").Value = 5396.4287
$ws.Range("I127").Value = 5396.4287
$ws.Range("K127").Value = 16189.2861
$ws.Range("M127").Value = -11229.2861
$ws.Range("H132").Value = 1383.78
$ws.Range("I132").Value = 1291.2554
$ws.Range("K132").Value = 3873.7662
$ws.Range("M132").Value = -1343.7662

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 62503660
$ws.Range("I2").Value = 1363.8889
$ws.Range("K2").Value = 1363.8889
$ws.Range("M2").Value = -1250.8889
$ws.Range("H61").Value = 6499.528
$ws.Range("J61").Value = 13483.333
$ws.Range("L61").Value = 13483.333
$ws.Range("N61").Value = -13907.333
$ws.Range("H88").Value = 2160.7
$ws.Range("J88").Value = 2315.2856
$ws.Range("L88").Value = 2315.2856
$ws.Range("N88").Value = -3127.2856
$ws.Range("H91").Value = 2160.7
$ws.Range("J91").Value = 2315.2856
$ws.Range("L91").Value = 2315.2856
$ws.Range("N91").Value = -5123.2856
$ws.Range("H116").Value = 62503660
$ws.Range("I116").Value = 1363.8889
$ws.Range("K116").Value = 1363.8889
$ws.Range("M116").Value = 930.1111000000001
$ws.Range("H122").Value = 3618.449
$ws.Range("I122").Value = 2735.6216
$ws.Range("J122").Value = 6340.5
$ws.Range("K122").Value = 8206.864799999999
$ws.Range("L122").Value = 19021.5
$ws.Range("M122").Value = -5756.864799999999
$ws.Range("N122").Value = -23921.5
$ws.Range("H132").Value = 3823.7273
$ws.Range("I132").Value = 1664.7872
$ws.Range("K132").Value = 4994.3616
$ws.Range("M132").Value = -2464.3616
$ws.Range("H136").Value = 6499.528
$ws.Range("J136").Value = 13483.333
$ws.Range("L136").Value = 40449.999
$ws.Range("N136").Value = -45549.999

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 62503660
$ws.Range("I3").Value = 1363.8889
$ws.Range("K3").Value = 1363.8889
$ws.Range("M3").Value = -1249.8889
$ws.Range("H20").Value = 5378307
$ws.Range("I20").Value = 9806128
$ws.Range("K20").Value = 9806128
$ws.Range("M20").Value = -9805881
$ws.Range("H86").Value = 62552572
$ws.Range("I86").Value = 14766851
$ws.Range("K86").Value = 14766851
$ws.Range("M86").Value = -14765728
$ws.Range("H89").Value = 62552572
$ws.Range("I89").Value = 14766851
$ws.Range("K89").Value = 73834255
$ws.Range("M89").Value = -73828639
$ws.Range("H94").Value = 1621
$ws.Range("I94").Value = 1376.125
$ws.Range("K94").Value = 1376.125
$ws.Range("M94").Value = -925.125
$ws.Range("H105").Value = 3398.1738
$ws.Range("I105").Value = 2594.8572
$ws.Range("K105").Value = 2594.8572
$ws.Range("M105").Value = -847.8571999999999
$ws.Range("H134").Value = 4315.9077
$ws.Range("I134").Value = 1745.7291
$ws.Range("K134").Value = 5237.1873
$ws.Range("M134").Value = -2702.1873

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5628
$ws.Range("I16").Value = 2895.875
$ws.Range("K16").Value = 2895.875
$ws.Range("M16").Value = -2608.875
$ws.Range("H58").Value = 8777378
$ws.Range("I58").Value = 16669132
$ws.Range("J58").Value = 8762.259
$ws.Range("K58").Value = 16669132
$ws.Range("L58").Value = 8762.259
$ws.Range("M58").Value = -16668929
$ws.Range("N58").Value = -9168.259
$ws.Range("H86").Value = 15634875
$ws.Range("I86").Value = 62500000
$ws.Range("J86").Value = 13166.667
$ws.Range("K86").Value = 62500000
$ws.Range("L86").Value = 13166.667
$ws.Range("M86").Value = -62498877
$ws.Range("N86").Value = -15412.667
$ws.Range("H88").Value = 229805.6
$ws.Range("I88").Value = 35000
$ws.Range("J88").Value = 278507
$ws.Range("K88").Value = 35000
$ws.Range("L88").Value = 278507
$ws.Range("M88").Value = -34594
$ws.Range("N88").Value = -279319
$ws.Range("H89").Value = 15634875
$ws.Range("I89").Value = 62500000
$ws.Range("J89").Value = 13166.667
$ws.Range("K89").Value = 312500000
$ws.Range("L89").Value = 65833.33499999999
$ws.Range("M89").Value = -312494384
$ws.Range("N89").Value = -77065.33499999999
$ws.Range("H91").Value = 229805.6
$ws.Range("I91").Value = 35000
$ws.Range("J91").Value = 278507
$ws.Range("K91").Value = 35000
$ws.Range("L91").Value = 278507
$ws.Range("M91").Value = -33596
$ws.Range("N91").Value = -281315
$ws.Range("H113").Value = 5628
$ws.Range("I113").Value = 2895.875
$ws.Range("K113").Value = 2895.875
$ws.Range("M113").Value = -725.875
$ws.Range("H122").Value = 2332.5715
$ws.Range("I122").Value = 1590.9166
$ws.Range("K122").Value = 4772.7498
$ws.Range("M122").Value = -2322.7498
$ws.Range("H132").Value = 8702023
$ws.Range("I132").Value = 4092.762
$ws.Range("J132").Value = 16008285
$ws.Range("K132").Value = 12278.286
$ws.Range("L132").Value = 48024855
$ws.Range("M132").Value = -9748.286
$ws.Range("N132").Value = -48029915
$ws.Range("H134").Value = 4369.0684
$ws.Range("I134").Value = 1795.6046
$ws.Range("K134").Value = 5386.8138
$ws.Range("M134").Value = -2851.8138
$ws.Range("H136").Value = 8777378
$ws.Range("I136").Value = 16669132
$ws.Range("J136").Value = 8762.259
$ws.Range("K136").Value = 50007396
$ws.Range("L136").Value = 26286.777
$ws.Range("M136").Value = -50004846
$ws.Range("N136").Value = -31386.777

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 70049.30499999999
$ws.Range("I2").Value = 12559.792
$ws.Range("J2").Value = 223354.67
$ws.Range("K2").Value = 75358.75199999999
$ws.Range("L2").Value = 1340128.02
$ws.Range("M2").Value = -75245.75199999999
$ws.Range("N2").Value = -1340354.02
$ws.Range("H5").Value = 4081.25
$ws.Range("I5").Value = 2400.1
$ws.Range("J5").Value = 6883.1665
$ws.Range("K5").Value = 7200.299999999999
$ws.Range("L5").Value = 20649.4995
$ws.Range("M5").Value = -7088.299999999999
$ws.Range("N5").Value = -20873.4995
$ws.Range("H121").Value = 10001642
$ws.Range("J121").Value = 2001.875
$ws.Range("L121").Value = 6005.625
$ws.Range("N121").Value = -8625.625
$ws.Range("H131").Value = 1920.5385
$ws.Range("I131").Value = 919.75
$ws.Range("J131").Value = 2220.775
$ws.Range("K131").Value = 2759.25
$ws.Range("L131").Value = 6662.325000000001
$ws.Range("M131").Value = 2280.75
$ws.Range("N131").Value = -16742.325
$ws.Range("H135").Value = 4081.25
$ws.Range("I135").Value = 2400.1
$ws.Range("J135").Value = 6883.1665
$ws.Range("K135").Value = 21600.9
$ws.Range("L135").Value = 61948.4985
$ws.Range("M135").Value = -19065.9
$ws.Range("N135").Value = -67018.4985

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 27034858
$ws.Range("I70").Value = 52637916
$ws.Range("K70").Value = 52637916
$ws.Range("M70").Value = -52637646
$ws.Range("H73").Value = 27034858
$ws.Range("I73").Value = 52637916
$ws.Range("K73").Value = 52637916
$ws.Range("M73").Value = -52636980
$ws.Range("H102").Value = 1577.7678
$ws.Range("I102").Value = 1394.2245
$ws.Range("K102").Value = 1394.2245
$ws.Range("M102").Value = 227.7755
$ws.Range("H122").Value = 1530143.1
$ws.Range("I122").Value = 2053377.6
$ws.Range("J122").Value = 4042.4167
$ws.Range("K122").Value = 6160132.800000001
$ws.Range("L122").Value = 12127.2501
$ws.Range("M122").Value = -6157682.800000001
$ws.Range("N122").Value = -17027.2501
$ws.Range("H132").Value = 4636.6
$ws.Range("I132").Value = 3381.848
$ws.Range("K132").Value = 10145.544
$ws.Range("M132").Value = -7615.544

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 201.6
$ws.Range("I16").Value = 212.69565
$ws.Range("J16").Value = 74
$ws.Range("K16").Value = 212.69565
$ws.Range("L16").Value = 74
$ws.Range("M16").Value = -42.69565
$ws.Range("N16").Value = -414
$ws.Range("H31").Value = 77475.53999999999
$ws.Range("I31").Value = 200545
$ws.Range("J31").Value = 557.125
$ws.Range("K31").Value = 200545
$ws.Range("L31").Value = 557.125
$ws.Range("M31").Value = -200297
$ws.Range("N31").Value = -1053.125

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3020.1667
$ws.Range("I126").Value = 1373.6666
$ws.Range("J126").Value = 4666.6665
$ws.Range("K126").Value = 4120.9998
$ws.Range("L126").Value = 13999.9995
$ws.Range("M126").Value = -1650.9998
$ws.Range("N126").Value = -18939.9995
